# Insert a new weekly price-report row for "Vega Monumental Concepción" /
# Arándano (blue) at row 54, pushing the existing rows 54-88 down to 55-89.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Shift rows 54:88 down by one row (creates a new blank row 54).
$ws.Rows("54:54").Insert()

# Populate the new row 54 with the latest weekly sample.
$ws.Range("A54").Value = 11
$ws.Range("B54").Value = "Vega Monumental Concepción"
$ws.Range("C54").Value = "Bíobío"
$ws.Range("D54").Value = 44596
$ws.Range("E54").Value = 8
$ws.Range("F54").Value = "Fruta"
$ws.Range("G54").Value = 100101
$ws.Range("H54").Value = "Berries"
$ws.Range("I54").Value = 100101001
$ws.Range("J54").Value = "Arándano (blue)"
$ws.Range("K54").Value = "Sin especificar"
$ws.Range("L54").Value = "Primera"
$ws.Range("M54").Value = 250
$ws.Range("N54").Value = 2800
$ws.Range("O54").Value = 3000
$ws.Range("P54").Value = 2920
$ws.Range("Q54").Value = "$/bandeja 2 kilos"
$ws.Range("R54").Value = "Provincia de Linares"
$ws.Range("S54").Value = 1460
$ws.Range("T54").Value = 2
